$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the required fields in the header row
$ws.Range("A1").Value = "Nombre (*obligatorio)"
$ws.Range("B1").Value = "Fecha de nacimiento (dd-mm-aa) (*obligatorio)"

# The longer header text needs a taller header row
$ws.Rows.Item(1).RowHeight = 46.25

# Move / extend the active selection to the header row
$ws.Range("A1:D1").Select()
